$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset gained a new weekly price record. It is inserted as row 57,
# pushing the former rows 57-74 down to become rows 58-75 (dimension grows
# from A1:R74 to A1:R75).
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new record's values.
$ws.Cells.Item(57, 1).Value = 3
$ws.Cells.Item(57, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(57, 3).Value = "Coquimbo"
$ws.Cells.Item(57, 4).Value = 44900
$ws.Cells.Item(57, 5).Value = 5
$ws.Cells.Item(57, 6).Value = 100112022
$ws.Cells.Item(57, 7).Value = "Arveja Verde"
$ws.Cells.Item(57, 8).Value = "Perfection"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 73
$ws.Cells.Item(57, 11).Value = 21000
$ws.Cells.Item(57, 12).Value = 22000
$ws.Cells.Item(57, 13).Value = 21479
$ws.Cells.Item(57, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(57, 15).Value = "Región Metropolitana"
$ws.Cells.Item(57, 16).Value = 859
$ws.Cells.Item(57, 17).Value = 25
$ws.Cells.Item(57, 18).Value = "Hortaliza"
